$wb = $excel.ActiveWorkbook

# Update the DDT_Rekomendasi Harga sheet values
$ws2 = $wb.Worksheets.Item("DDT_Rekomendasi Harga")
$ws2.Range("D3").Value = 5000000
$ws2.Range("D4").Value = 2000000

# Make DDT_Rekomendasi Harga the active/selected sheet and set its selection
$ws2.Activate()
$ws2.Range("D6").Select()
